$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the Kyrgyz title (A1, shared-string backed) - wording/terminology fix
$ws.Range("A1").Value = "8.10.2.2 Камсыздандыруу компаниялардын финансылык көрсөткүчтөрү"

# Append the new 2023 reporting column (Q), mirroring the existing 2022 column (P)
$ws.Range("Q3").Value = 2023
$ws.Range("Q4").Value = 16
$ws.Range("Q5").Value = 3031.4

# Copy the formatting (styles/number formats) of column P onto the new column Q
$ws.Range("P3:P5").Copy()
$ws.Range("Q3:Q5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Return the cursor to the top-left cell
$ws.Range("A1").Select()
